# Auto-generated edit script applying cryptos.xlsx diff (updates rankings/prices/volume columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure all target cells are formatted as text so numeric-looking strings
# (e.g. "241.89", "0.0747") are preserved verbatim instead of being parsed as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.437.84"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.37%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.052.63"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.89"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.38%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.27"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -5.93%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "58.26"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.64%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0747"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.58%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.902"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.65"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -5.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.354.85"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.37"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -5.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.051.22"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.378.41"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.73"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -7.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.77"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0853"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "238.28"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.22"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.21%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.36"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.11"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "164.14"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.98"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.53%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.30%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +8.16%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -6.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.43"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.94%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.34%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.83"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.18"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0816"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -5.64%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -5.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.83"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.78%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.21%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Cronos"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0932"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -5.12%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.82"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -9.22%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.10"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.32%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "93.54"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.399.68"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +8.54%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.54"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +12.39%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.84"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -6.59%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.242.58"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.86%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.24"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.77%  "
